# 11 May 2024 Selenium Data Driven
#
# The "TestData" sheet drives a Selenium data-driven test: column A
# ("RunMode") flags which rows should execute ("Yes"/"No") and column S
# ("Result") records the outcome of the latest run ("Pass"/"Fail").
#
# This edit:
#   1) fixes row 4's RunMode, which was mistakenly left as "No" (it should
#      run, like the other rows), and
#   2) records the results of the latest Selenium run: every row flagged
#      to run ("Yes") passed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Row 4 should participate in the run, not be skipped.
$ws.Cells.Item(4, 1).Value = "Yes"

# Stamp the Result column (S) with the outcome of the data-driven run:
# every row whose RunMode is "Yes" passed.
for ($r = 2; $r -le 8; $r++) {
    $runMode = $ws.Cells.Item($r, 1).Value2
    if ($runMode -eq "Yes") {
        $ws.Cells.Item($r, 19).Value = "Pass"
    }
}

# Leave the selection where the author left it after the run.
$ws.Range("A4").Select()
